$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '68.356.32'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -2.91%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.664.21'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -4.29%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '591.62'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.07%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '180.11'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +8.19%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '3.659.45'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -4.22%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.629'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -4.94%  '
$ws.Range('E9').Value = '  +0.39%  '
$ws.Range('E10').Value = '  -3.50%  '
$ws.Range('E11').Value = '  -7.34%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '56.27'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +7.13%  '
$ws.Range('E13').Value = '  -8.14%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '10.63'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -5.46%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '4.248.18'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -4.16%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.646.54'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -4.72%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '19.30'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -6.56%  '
$ws.Range('E18').Value = '  -1.97%  '
$ws.Range('E19').Value = '  -6.70%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '1.12'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -6.37%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '68.058.95'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -3.11%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '409.36'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -5.55%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.58'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -3.27%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '88.50'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -4.80%  '
$ws.Range('B25').Value = 'ImmutableX'
$ws.Range('C25').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '3.02'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -7.44%  '
$ws.Range('B26').Value = 'InternetComputer(DFINITY)'
$ws.Range('C26').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '12.74'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -6.89%  '
$ws.Range('B27').Value = 'RenderToken'
$ws.Range('C27').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.86'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -5.04%  '
$ws.Range('B28').Value = 'Toncoin'
$ws.Range('C28').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '3.93'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -1.27%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '6.05'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +1.49%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '9.48'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -8.89%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '32.60'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -6.15%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '7.22'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -11.72%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '12.37'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -7.61%  '
$ws.Range('E34').Value = '  -5.57%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '64.69'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -5.23%  '
$ws.Range('B36').Value = 'Bittensor'
$ws.Range('C36').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '604.22'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -5.21%  '
$ws.Range('B37').Value = 'InjectiveProtocol'
$ws.Range('C37').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '43.23'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -9.56%  '
$ws.Range('E38').Value = '  -9.72%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.999'
$ws.Range('D39').Style = 'Normal'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.399'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -6.72%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.999'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.03%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.136'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -6.08%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '3.03'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -5.62%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.71'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -7.16%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0437'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -6.07%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.87'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -11.81%  '
$ws.Range('B47').Value = 'WEMIXToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.72'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -1.64%  '
$ws.Range('B48').Value = 'Stellar'
$ws.Range('C48').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.135'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -5.42%  '
$ws.Range('B49').Value = 'THORChain'
$ws.Range('C49').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '9.02'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -8.41%  '
$ws.Range('B50').Value = 'ApeXProtocol'
$ws.Range('C50').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '3.13'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -3.81%  '
$ws.Range('B51').Value = 'Maker'
$ws.Range('C51').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.713.74'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -5.23%  '
